# Apply "update database and change read_price algorithm":
# - drop the oldest quarter column (column D) and shift every later quarter left
# - append the newest quarter (column M) with its header, publish-date and data
# - one previously-published date/revision label gets amended in place

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Delete the oldest quarter column; everything to the right shifts left one column.
$ws.Range("D:D").Delete()

# 2) Make sure the (now last) data column M keeps the same 31-width style as the
#    other "publish date" columns (E, I, M were width 31 before the shift).
$ws.Range("M1").ColumnWidth = 31

# 3) Populate the newly shifted-in last column (M) with the newest quarter's data.
$ws.Range("M8").Value  = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value  = "1402-02-28"
$ws.Range("M11").Value = 12501
$ws.Range("M12").Value = -10156
$ws.Range("M13").Value = 2345
$ws.Range("M14").Value = -522
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 289
$ws.Range("M17").Value = 2112
$ws.Range("M18").Value = -398
$ws.Range("M19").Value = -143
$ws.Range("M20").Value = 1571
$ws.Range("M21").Value = -79
$ws.Range("M22").Value = 1492
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 1492
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 28601
$ws.Range("M27").Value = 0

# 4) The previously-published "1401-10-28 (6)" label (now shifted into column I)
#    was amended by the issuer to "1402-02-28 (7)".
$ws.Range("I9").Value = "1402-02-28 (7)"
